# Update Excel file with latest predictions

function Set-RowValues($ws, $r, $vals) {
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
    }
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Home win": replace row 2, drop row 3 (Guatemala match removed)
# ---------------------------------------------------------------------------
$wsHome = $wb.Worksheets.Item("Home win")
$wsHome.Rows.Item(3).Delete()
Set-RowValues $wsHome 2 @("04-02-2025 20:45", "SCOTLAND", "LEAGUE TWO", "Elgin City - Bonnyrigg Rose Athletic", 73.3, 1.91)

# ---------------------------------------------------------------------------
# Sheet "Away Win": row 2 stays as-is, drop row 3 (Uruguay match removed)
# ---------------------------------------------------------------------------
$wsAway = $wb.Worksheets.Item("Away Win")
$wsAway.Rows.Item(3).Delete()

# ---------------------------------------------------------------------------
# Sheet "Draw": replace rows 2-3 with new matches, drop old rows 4-6
# ---------------------------------------------------------------------------
$wsDraw = $wb.Worksheets.Item("Draw")
$wsDraw.Rows.Item(6).Delete()
$wsDraw.Rows.Item(5).Delete()
$wsDraw.Rows.Item(4).Delete()
Set-RowValues $wsDraw 2 @("04-02-2025 20:45", "ENGLAND", "CHAMPIONSHIP", "Burnley - Oxford United", 70, 4.5)
Set-RowValues $wsDraw 3 @("07-01-2025 20:00", "ENGLAND", "EFL TROPHY", "Port Vale - Wrexham", 73.3, 3.4)

# ---------------------------------------------------------------------------
# Sheet "Btts": replace rows 2-7, add new rows 8-9
# ---------------------------------------------------------------------------
$wsBtts = $wb.Worksheets.Item("Btts")
Set-RowValues $wsBtts 2 @("03-02-2025 18:00", "CZECH-REPUBLIC", "CZECH LIGA", "Sigma Olomouc - Plzen", 76.7, 1.75)
Set-RowValues $wsBtts 3 @("03-02-2025 13:30", "EGYPT", "SECOND LEAGUE", "Dayrout - Kahraba Ismailia", 76, 2)
Set-RowValues $wsBtts 4 @("03-02-2025 13:30", "EGYPT", "SECOND LEAGUE", "Raya Ghazl - La Viena FC", 80, 2.1)
Set-RowValues $wsBtts 5 @("03-02-2025 14:10", "SAUDI-ARABIA", "DIVISION 1", "Ohod - Abha", 88, 1.7)
Set-RowValues $wsBtts 6 @("23-11-2024 16:00", "ENGLAND", "LEAGUE TWO", "Salford City - Bromley", 76.7, 1.83)
Set-RowValues $wsBtts 7 @("05-02-2025 00:00", "CHILE", "COPA CHILE", "Deportes Limache - Union San Felipe", 83.3, 1.73)
Set-RowValues $wsBtts 8 @("04-02-2025 21:10", "FRANCE", "COUPE DE FRANCE", "Le Mans - Paris Saint Germain", 90, 2.1)
Set-RowValues $wsBtts 9 @("04-02-2025 19:00", "FRANCE", "COUPE DE FRANCE", "Lille - Dunkerque", 78.3, 1.85)

# ---------------------------------------------------------------------------
# Sheet "Over_Under": replace rows 2-5, add new rows 6-8
# ---------------------------------------------------------------------------
$wsOU = $wb.Worksheets.Item("Over_Under")
Set-RowValues $wsOU 2 @("03-02-2025 21:45", "PORTUGAL", "PRIMEIRA LIGA", "Rio Ave - FC Porto", 80, 1.75, 45, 2.75)
Set-RowValues $wsOU 3 @("03-02-2025 14:10", "SAUDI-ARABIA", "DIVISION 1", "Ohod - Abha", 80, 1.93, 46.7, 3.35)
Set-RowValues $wsOU 4 @("04-02-2025 20:45", "ENGLAND", "EFL TROPHY", "Stevenage - Birmingham", 80, 1.91, 13.3, 3.1)
Set-RowValues $wsOU 5 @("04-02-2025 20:45", "ENGLAND", "NATIONAL LEAGUE - NORTH", "Chorley - Buxton", 85, 1.75, 50, 3)
Set-RowValues $wsOU 6 @("04-02-2025 20:45", "ENGLAND", "NATIONAL LEAGUE - NORTH", "Scarborough Athletic - Radcliffe", 70, 1.65, 60, 2.6)
Set-RowValues $wsOU 7 @("04-02-2025 20:45", "ENGLAND", "NON LEAGUE PREMIER - ISTHMIAN", "Cray Valley PM - Lewes", 73.3, 1.57, 60, 2.4)
Set-RowValues $wsOU 8 @("04-02-2025 19:00", "FRANCE", "COUPE DE FRANCE", "Lille - Dunkerque", 70, 1.8, 60, 3)
